# Auto-generated edit script applying the Yojimbo_Profits value updates.
# For each affected Leve row, sets the recomputed market-price/profit
# columns (H..N) per the commit diff, and clears cells that the diff
# removes entirely (kept blank, matching the source row layout).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!100 - Asking for a Friend (Beetle Glue)
$ws.Range("H100").Value = 4099.875
$ws.Range("I100").Value = 3199.75
$ws.Range("K100").Value = 3199.75
$ws.Range("M100").Value = -2658.75

# ALC!127 - Liquid Competence (Competent Craftsman's Draught)
$ws.Range("H127").Value = 2002.5
$ws.Range("I127").Value = 379.5
$ws.Range("J127").Value = 2543.5
$ws.Range("K127").Value = 1138.5
$ws.Range("L127").Value = 7630.5
$ws.Range("M127").Value = 3821.5
$ws.Range("N127").Value = -17550.5

# ALC!137 - Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws.Range("H137").Value = 2478.6545
$ws.Range("I137").Value = 2284.6365
$ws.Range("J137").Value = 3254.7273
$ws.Range("K137").Value = 6853.9095
$ws.Range("L137").Value = 9764.1819
$ws.Range("M137").Value = -4303.9095
$ws.Range("N137").Value = -14864.1819

$ws = $wb.Worksheets.Item("ARM")
# ARM!2 - Ain't Got No Ingots (Bronze Ingot)
$ws.Range("H2").Value = 1225.8096
$ws.Range("I2").Value = 1145.3684
$ws.Range("J2").Value = 1990
$ws.Range("K2").Value = 1145.3684
$ws.Range("L2").Value = 1990
$ws.Range("M2").Value = -1032.3684
$ws.Range("N2").Value = -2216

# ARM!74 - As the Bolt Flies (Titanium Nugget)
$ws.Range("H74").Value = 2821.149
$ws.Range("I74").Value = 3362.5557
$ws.Range("K74").Value = 3362.5557
$ws.Range("M74").Value = -2488.5557

# ARM!77 - Heavy Metal Banned (L) (Titanium Nugget)
$ws.Range("H77").Value = 2821.149
$ws.Range("I77").Value = 3362.5557
$ws.Range("K77").Value = 16812.7785
$ws.Range("M77").Value = -12444.7785

# ARM!97 - Ore for Me (High Steel Ingot)
$ws.Range("H97").Value = 1737.9667
$ws.Range("I97").Value = 792.04346
$ws.Range("J97").Value = 4846
$ws.Range("K97").Value = 792.04346
$ws.Range("L97").Value = 4846
$ws.Range("M97").Value = -296.04346
$ws.Range("N97").Value = -5838

# ARM!116 - No Scope (Titanbronze Ingot)
$ws.Range("H116").Value = 1225.8096
$ws.Range("I116").Value = 1145.3684
$ws.Range("J116").Value = 1990
$ws.Range("K116").Value = 1145.3684
$ws.Range("L116").Value = 1990
$ws.Range("M116").Value = 1148.6316
$ws.Range("N116").Value = -6578

# ARM!122 - Haste for High Durium (High Durium Nugget)
$ws.Range("H122").Value = 1462
$ws.Range("I122").Value = 1260.7826
$ws.Range("J122").Value = 2233.3333
$ws.Range("K122").Value = 3782.3478
$ws.Range("L122").Value = 6699.999899999999
$ws.Range("M122").Value = -1332.3478
$ws.Range("N122").Value = -11599.9999

# ARM!132 - Don't Bore Me, Ore Me (Mountain Chromite Ingot)
$ws.Range("H132").Value = 6909.6484
$ws.Range("I132").Value = 4792.107
$ws.Range("J132").Value = 13497.556
$ws.Range("K132").Value = 14376.321
$ws.Range("L132").Value = 40492.66800000001
$ws.Range("M132").Value = -11846.321
$ws.Range("N132").Value = -45552.66800000001

# ARM!133 - Shielding My Students (Mountain Chromite Tower Shield)
$ws.Range("H133").Value = 26972.875
$ws.Range("J133").Value = 26972.875
$ws.Range("L133").Value = 26972.875
$ws.Range("N133").Value = -32032.875

$ws = $wb.Worksheets.Item("BSM")
# BSM!3 - Hells Bells (Bronze Ingot)
$ws.Range("H3").Value = 1225.8096
$ws.Range("I3").Value = 1145.3684
$ws.Range("J3").Value = 1990
$ws.Range("K3").Value = 1145.3684
$ws.Range("L3").Value = 1990
$ws.Range("M3").Value = -1031.3684
$ws.Range("N3").Value = -2218

# BSM!105 - Ingot to Wing It (Molybdenum Ingot)
$ws.Range("H105").Value = 71432456
$ws.Range("I105").Value = 100003736
$ws.Range("J105").Value = 4249.75
$ws.Range("K105").Value = 100003736
$ws.Range("L105").Value = 4249.75
$ws.Range("M105").Value = -100001989
$ws.Range("N105").Value = -7743.75

# BSM!134 - Ruthenium Supremium (Ruthenium Ingot)
$ws.Range("H134").Value = 4639.6206
$ws.Range("I134").Value = 5133.4346
$ws.Range("J134").Value = 2746.6667
$ws.Range("K134").Value = 15400.3038
$ws.Range("L134").Value = 8240.000100000001
$ws.Range("M134").Value = -12865.3038
$ws.Range("N134").Value = -13310.0001

$ws = $wb.Worksheets.Item("CRP")
# CRP!22 - Driving Up the Wall (Elm Lumber)
$ws.Range("H22").Value = 750.1905
$ws.Range("I22").Value = 861.375
$ws.Range("J22").Value = 394.4
$ws.Range("K22").Value = 861.375
$ws.Range("L22").Value = 394.4
$ws.Range("M22").Value = -511.375
$ws.Range("N22").Value = -1094.4

# CRP!62 - Splinter in the Sewers (Cedar Lumber)
$ws.Range("H62").Value = 3236.6
$ws.Range("I62").Value = 2317.5
$ws.Range("J62").Value = 4615.25
$ws.Range("K62").Value = 2317.5
$ws.Range("L62").Value = 4615.25
$ws.Range("M62").Value = -1693.5
$ws.Range("N62").Value = -5863.25

# CRP!65 - The Lumber of Their Discontent (L) (Cedar Lumber)
$ws.Range("H65").Value = 3236.6
$ws.Range("I65").Value = 2317.5
$ws.Range("J65").Value = 4615.25
$ws.Range("K65").Value = 11587.5
$ws.Range("L65").Value = 23076.25
$ws.Range("M65").Value = -8467.5
$ws.Range("N65").Value = -29316.25

# CRP!74 - License to Heal (Dark Chestnut Rod)
$ws.Range("H74").Value = 40267.273
$ws.Range("J74").Value = 40267.273
$ws.Range("L74").Value = 40267.273
$ws.Range("N74").Value = -42015.273

# CRP!77 - Purified Polyrhythm (L) (Dark Chestnut Rod)
$ws.Range("H77").Value = 40267.273
$ws.Range("J77").Value = 40267.273
$ws.Range("L77").Value = 120801.819
$ws.Range("N77").Value = -129537.819

$ws = $wb.Worksheets.Item("CUL")
# CUL!107 - Slippery Service (Frantoio Oil)
$ws.Range("H107").Value = 579.4545000000001
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 579.4545000000001
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1738.3635
$ws.Range("N107").Value = -5578.3635
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# GSM!102 - Put the Metal to the Peddle (Durium Ingot)
$ws.Range("H102").Value = 1296.4
$ws.Range("I102").Value = 1082.6842
$ws.Range("K102").Value = 1082.6842
$ws.Range("M102").Value = 539.3158000000001

# GSM!126 - Gold Rush Order (Phrygian Gold Ingot)
$ws.Range("H126").Value = 1169.4546
$ws.Range("I126").Value = 1100.5714
$ws.Range("J126").Value = 1290
$ws.Range("K126").Value = 3301.7142
$ws.Range("L126").Value = 3870
$ws.Range("M126").Value = -831.7142000000003
$ws.Range("N126").Value = -8810

# GSM!132 - On Board for Lar (Lar Ingot)
$ws.Range("H132").Value = 5315.263
$ws.Range("I132").Value = 4007.756
$ws.Range("J132").Value = 8665.75
$ws.Range("K132").Value = 12023.268
$ws.Range("L132").Value = 25997.25
$ws.Range("M132").Value = -9493.268
$ws.Range("N132").Value = -31057.25

$ws = $wb.Worksheets.Item("LTW")
# LTW!7 - Tan Before the Ban (Leather)
$ws.Range("H7").Value = 1444.8889
$ws.Range("I7").Value = 1125.5
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 1125.5
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -1013.5
$ws.Range("N7").Value = -4224

# LTW!16 - Saddle Sore (Hard Leather)
$ws.Range("H16").Value = 902.1
$ws.Range("I16").Value = 805
$ws.Range("J16").Value = 1776
$ws.Range("K16").Value = 805
$ws.Range("L16").Value = 1776
$ws.Range("M16").Value = -635
$ws.Range("N16").Value = -2116

# LTW!122 - Hell on Leather (Gaja Leather)
$ws.Range("H122").Value = 3699.4849
$ws.Range("I122").Value = 3975.3635
$ws.Range("J122").Value = 3147.7273
$ws.Range("K122").Value = 11926.0905
$ws.Range("L122").Value = 9443.1819
$ws.Range("M122").Value = -9476.0905
$ws.Range("N122").Value = -14343.1819

# LTW!126 - Battered Books (Saiga Leather)
$ws.Range("H126").Value = 1444.8889
$ws.Range("I126").Value = 1125.5
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 3376.5
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -906.5
$ws.Range("N126").Value = -16940

# LTW!132 - Tenets of Tanning (Silver Lobo Leather)
$ws.Range("H132").Value = 3901.966
$ws.Range("I132").Value = 4251.2324
$ws.Range("J132").Value = 2963.3125
$ws.Range("K132").Value = 12753.6972
$ws.Range("L132").Value = 8889.9375
$ws.Range("M132").Value = -10223.6972
$ws.Range("N132").Value = -13949.9375

# LTW!136 - Respect for Br'aax (Br'aax Leather)
$ws.Range("H136").Value = 2885.1777
$ws.Range("I136").Value = 2386.9546
$ws.Range("J136").Value = 3361.739
$ws.Range("K136").Value = 7160.8638
$ws.Range("L136").Value = 10085.217
$ws.Range("M136").Value = -4610.8638
$ws.Range("N136").Value = -15185.217

$ws = $wb.Worksheets.Item("WVR")
# WVR!76 - Finger on the Pulse (Ramie Halfgloves of Healing)
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# WVR!79 - Chirurgeon Hand in Glove (L) (Ramie Halfgloves of Healing)
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# WVR!113 - A Tender Table (Pixie Floss)
$ws.Range("H113").Value = 340.33334
$ws.Range("I113").Value = 245
$ws.Range("J113").Value = 388
$ws.Range("K113").Value = 735
$ws.Range("L113").Value = 1164
$ws.Range("M113").Value = 1435
$ws.Range("N113").Value = -5504

# WVR!125 - Color Coated (Almasty Serge Coat of Healing)
$ws.Range("H125").Value = 39729.062
$ws.Range("J125").Value = 39729.062
$ws.Range("L125").Value = 39729.062
$ws.Range("N125").Value = -49569.062
